# notas_ia.xlsx -- "fechando o repositorio com as notas do semestre 2025.1"
#
# Prova2 (sheet2): fix the H column averages so every row reads its own
# D/F cells (several had been copy-pasted from other rows), correct a
# handful of "Faltou na prova1?" (D) flags, drop leftover/empty C cells,
# normalize the E column format, de-duplicate a few redundant direct
# cell styles, and reset the sheet's saved scroll position.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Prova2")

# ---------------------------------------------------------------------
# 1) Header row (row 1): A1/C1 and D1/E1 were carrying their own
#    one-off style entries that duplicate existing ones (s=13 == s=9,
#    s=14 == s=1/0). Re-point them at the de-duplicated styles by
#    pasting formats from cells that already carry the target look.
# ---------------------------------------------------------------------
$ws2.Range("A2").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("C1").PasteSpecial(-4122)

$ws1 = $wb.Worksheets.Item("Prova1")
$ws1.Range("F1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("E1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Column E (rows 2-35): normalize the "DP" style down to the plain
#    numeric style already used by column D (s=10) instead of its own
#    duplicate (s=15).
# ---------------------------------------------------------------------
$ws2.Range("D2").Copy()
for ($r = 2; $r -le 35; $r++) {
    $ws2.Range("E$r").PasteSpecial(-4122)
}

# Rows 36-38: E had stray "22" (Regime Disciplinar) text values with yet
# another duplicate style (s=17). Clear them back to blank D10-style
# cells like the rest of the column.
foreach ($r in @(36, 37, 38)) {
    $ws2.Range("E$r").ClearContents()
    $ws2.Range("D2").Copy()
    $ws2.Range("E$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 3) Column C: drop the leftover empty/placeholder cells entirely so
#    the row goes back to just A/B/D/E/F/H.
# ---------------------------------------------------------------------
foreach ($r in @(5, 6, 10, 16, 18, 21, 23, 36, 37, 38)) {
    $ws2.Range("C$r").ClearContents()
}

# ---------------------------------------------------------------------
# 4) Column D: a handful of "Faltou na prova1?" flags were left at 0
#    when they should record how many faltas; fix the values and align
#    their style with the one other non-default D style in the sheet.
# ---------------------------------------------------------------------
$ws2.Range("D29").Value = 1
$ws2.Range("D30").Value = 1
$ws2.Range("D33").Value = 4
$ws2.Range("D34").Value = 3
$ws2.Range("D35").Value = 3

$ws2.Range("H29").Copy()
foreach ($r in @(29, 30, 33, 34, 35)) {
    $ws2.Range("D$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 5) Column H: every row should average its OWN D/F cells. Several rows
#    had been left pointing at a different row's D/F after past
#    sorting/copying -- fix the formulas first (this also refreshes the
#    cached <v> through recalculation).
# ---------------------------------------------------------------------
$hFixRows = @(5, 8, 10, 11, 14, 17, 21, 23, 24, 30, 31, 37)
foreach ($r in $hFixRows) {
    $ws2.Range("H$r").Formula = "=(D$r+F$r)/2"
}

# All of column H (rows 2-38) also carried its own direct style
# (s=10, s=5 or s=16) that isn't needed -- strip it back down to the
# workbook default, matching plain cells like F2.
$ws2.Range("F2").Copy()
for ($r = 2; $r -le 38; $r++) {
    $ws2.Range("H$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 6) Reset the saved scroll position for this sheet (it had been left
#    scrolled to row 16).
# ---------------------------------------------------------------------
$ws2.Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("A1").Select()

$wb.Save()
